$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48, shifting existing rows 48:120 down to 49:121
$ws.Rows(48).Insert()

# Populate the newly inserted row 48 with the new record
$ws.Range("A48").Value = 1
$ws.Range("B48").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C48").Value = "Arica y Parinacota"
$ws.Range("D48").Value = 45117
$ws.Range("E48").Value = 15
$ws.Range("F48").Value = 100112040
$ws.Range("G48").Value = "Cilantro"
$ws.Range("H48").Value = "Sin especificar"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 420
$ws.Range("K48").Value = 800
$ws.Range("L48").Value = 1000
$ws.Range("M48").Value = 871
$ws.Range("N48").Value = "`$/atado 1,5 a 2 kilos"
$ws.Range("O48").Value = "Región de Arica y Parinacota"
$ws.Range("P48").Value = 436
$ws.Range("Q48").Value = 2
$ws.Range("R48").Value = "Hortaliza"
